$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the whole "License Information" (Heading2) paragraph.
# ---------------------------------------------------------------------------
$d.Paragraphs(4).Range.Delete()

# ---------------------------------------------------------------------------
# 2. Remove the whole "This PDF version is provided under the same
#    license." paragraph (it is now paragraph 5, since paragraph 4 above
#    was just removed).
# ---------------------------------------------------------------------------
$d.Paragraphs(5).Range.Delete()

# ---------------------------------------------------------------------------
# 3. Rewrite the remaining license paragraph (now paragraph 4).
# ---------------------------------------------------------------------------
$p = $d.Paragraphs(4)

# Locate the block of text/hyperlinks that must be removed: everything from
# right after the bold "أسئلة الترجمة (unfoldingWord)" run up to and
# including the closing "CC BY-SA 4.0 license." (this also removes the two
# hyperlink runs in between).
$startRng = $p.Range.Duplicate()
$startRng.Find.Execute(" (Arabic) is based on", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$startPos = $startRng.Start

$endRng = $p.Range.Duplicate()
$endRng.Find.Execute("CC BY-SA 4.0 license.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$endPos = $endRng.End

$d.Range($startPos, $endPos).Text = ""

# Change the bold run's text from the Arabic title to the new English title.
$boldRng = $p.Range.Duplicate()
$boldRng.Find.Execute("أسئلة الترجمة (unfoldingWord)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$boldRng.Text = "unfoldingWord® Translation Questions"

# Replace the leftover "." with the full new license / adaptation text.
$periodRng = $p.Range.Duplicate()
$periodRng.Find.Execute(".", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$newText = " © 2022 unfoldingWord. Released under CC BY-SA 4.0 license. " + `
    "unfoldingWord® Translation Questions" + `
    " has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文) from " + `
    "unfoldingWord® Translation Questions" + `
    " © 2022 unfoldingWord. Released under CC BY-SA 4.0 license by Mission Mutual"

$periodRng.Text = $newText
